# feat: add 2022-Q4 data
#
# The workbook has a "总计" (totals) summary sheet followed by one sheet
# per quarter (2022-Q3, 2022-Q2, 2022-Q1). This script inserts a new
# "2022-Q4" quarter sheet right after "总计" (pushing the older quarters
# down, unchanged), fills it in with that quarter's fund data, and adds
# the corresponding summary row on "总计".

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")
$q3Sheet = $wb.Worksheets.Item("2022-Q3")

# --- 1. Create the new "2022-Q4" sheet -------------------------------------
# Copy the "2022-Q3" sheet (same layout/styles/page setup) and drop it
# right after "总计", then rename it and swap in the Q4 numbers.
$q3Sheet.Copy($null, $totalSheet)
$q4Sheet = $wb.Worksheets.Item(2)
$q4Sheet.Name = "2022-Q4"

$q4Sheet.Range("D2").Value = "'0.67"
$q4Sheet.Range("E2").Value = "'91.81"
$q4Sheet.Range("F2").Value = "'4.73"
$q4Sheet.Range("G2").Value = "'0.0317"
$q4Sheet.Range("H2").Value = 6

# --- 2. Update the "总计" summary sheet -------------------------------------
# Insert a fresh row for 2022-Q4 right under the header, reuse the
# formatting of the existing index column for the row that now falls off
# the bottom (2022-Q1), and rewrite the whole table with the new figures.
$totalSheet.Range("A2").Copy()
$totalSheet.Range("A5").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.03

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q3"
$totalSheet.Range("C3").Value = 1
$totalSheet.Range("D3").Value = 0.03

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2022-Q2"
$totalSheet.Range("C4").Value = 1
$totalSheet.Range("D4").Value = 0.03

$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = "2022-Q1"
$totalSheet.Range("C5").Value = 1
$totalSheet.Range("D5").Value = 0.02
